$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data Model")

# The "Customer" table (row 2 header) loses its "CustomerID" column; the
# remaining columns shift one to the left (B:L instead of C:M), and the
# trailing cell (M2, formerly CVV) is cleared since the range shrinks by one.
$ws.Range("B2").Value = "EmailID"
$ws.Range("C2").Value = "FullName"
$ws.Range("D2").Value = "Address"
$ws.Range("E2").Value = "City"
$ws.Range("F2").Value = "State"
$ws.Range("G2").Value = "Zip"
$ws.Range("H2").Value = "MobNo"
$ws.Range("I2").Value = "CreditCardNo"
$ws.Range("J2").Value = "CreditCardType"
$ws.Range("K2").Value = "ExpDate"
$ws.Range("L2").Value = "CVV"
$ws.Range("M2").ClearContents()

# The "Order" table's foreign-key column (H5) now references EmailID
# instead of the removed CustomerID column.
$ws.Range("H5").Value = "EmailID"

# Leave the selection where the edit happened, matching the saved UI state.
$ws.Range("H5").Select()
